# New clock configuration files
# Insert a new configuration row ("R1Av0004") above the existing "R1Bv0001"
# row (old row 16), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift row 15 (and everything below) down one row, inserting a fresh blank
# row 15 for the new configuration entry.
$ws.Rows("15:15").Insert() | Out-Null

# The insert carries a formatting-only placeholder cell down from row 14
# (column E) that has no counterpart in the new row's data; drop it.
$ws.Range("E15").Clear() | Out-Null

# Populate the new "R1Av0004" configuration row.
$ws.Range("A15").Value = "R1Av0004"
$ws.Range("C15").Value = 40
$ws.Range("H15:R15").Value = 320

# Match the saved selection/active cell.
$ws.Range("E15").Select() | Out-Null
